$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The input data format changes from:
#   name | pan | street_num | street_name | house_num | locality | city | state | pin | portalpass
# to:
#   name | phone | email | pan | aadhaar | type | portalpass | street_num | street_name | house_num | locality | city | state | pin
#
# i.e. two new columns (phone, email) are inserted right after "name", two more
# new columns (aadhaar, type) are inserted right after "pan", and the existing
# "portalpass" column is relocated to sit right after the new "type" column.

# 1) Insert two blank columns for "phone" and "email" right after "name" (old col A).
#    This pushes pan..portalpass (old B..J) right by two (pan -> D, portalpass -> L).
$ws.Columns("B:C").Insert()

# 2) Insert two blank columns for "aadhaar" and "type" right after "pan" (now col D).
#    This pushes street_num..portalpass right by two more (street_num -> G .. portalpass -> N).
$ws.Columns("E:F").Insert()

# 3) Relocate "portalpass" (now col N) so that it sits right after "type" (col F),
#    i.e. right before "street_num" (col G).
$ws.Columns("N:N").Cut()
$ws.Columns("G:G").Insert()

# The cut/insert leaves stray width markers behind on the now-unused trailing
# columns (L:N) - drop that leftover column formatting so they fall back to
# the sheet's default column width.
$ws.Columns("L:N").ClearFormats()

# 4) Fill in the headers for the newly inserted columns.
$ws.Range("B1").Value = "phone"
$ws.Range("C1").Value = "email"
$ws.Range("E1").Value = "aadhaar"
$ws.Range("F1").Value = "type"

# Match the saved selection state.
$ws.Range("F4").Select()
